$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matching the source data's inlineStr cell
# type) even when the string looks like a plain number (e.g. "12.42"),
# without leaving a stray NumberFormat/quotePrefix style behind on the cell.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" '67.259.52'
Set-TextValue "E2" '  +7.18%  '
Set-TextValue "D3" '3.597.46'
Set-TextValue "E3" '  +3.66%  '
Set-TextValue "E4" '  +0.08%  '
Set-TextValue "D5" '415.38'
Set-TextValue "E5" '  +0.33%  '
Set-TextValue "D6" '129.72'
Set-TextValue "E6" '  -0.58%  '
Set-TextValue "D7" '0.652'
Set-TextValue "E7" '  +3.79%  '
Set-TextValue "D8" '3.589.39'
Set-TextValue "E8" '  +3.63%  '
Set-TextValue "E9" '  -0.03%  '
Set-TextValue "D10" '0.781'
Set-TextValue "E10" '  +7.28%  '
Set-TextValue "D11" '0.177'
Set-TextValue "E11" '  +17.15%  '
Set-TextValue "D12" '0.0000338'
Set-TextValue "E12" '  +53.70%  '
Set-TextValue "D13" '42.49'
Set-TextValue "E13" '  -0.37%  '
Set-TextValue "D14" '9.91'
Set-TextValue "E14" '  +3.00%  '
Set-TextValue "D15" '4.169.26'
Set-TextValue "E15" '  +3.74%  '
Set-TextValue "E16" '  -0.32%  '
Set-TextValue "D17" '20.31'
Set-TextValue "E17" '  -1.30%  '
Set-TextValue "D18" '3.600.93'
Set-TextValue "E18" '  +3.82%  '
Set-TextValue "E19" '  +5.34%  '
Set-TextValue "D20" '67.118.78'
Set-TextValue "E20" '  +6.99%  '
Set-TextValue "D21" '12.28'
Set-TextValue "E21" '  -2.91%  '
Set-TextValue "D22" '452.04'
Set-TextValue "E22" '  -2.49%  '
Set-TextValue "D23" '89.52'
Set-TextValue "E23" '  -1.37%  '
Set-TextValue "E24" '  -3.66%  '
Set-TextValue "D25" '13.12'
Set-TextValue "E25" '  -1.42%  '
Set-TextValue "E26" '  +0.89%  '
Set-TextValue "D27" '9.97'
Set-TextValue "E27" '  -7.30%  '
Set-TextValue "D28" '35.30'
Set-TextValue "E28" '  +5.34%  '
Set-TextValue "D29" '4.87'
Set-TextValue "E29" '  +1.47%  '
Set-TextValue "D30" '12.42'
Set-TextValue "E31" '  +3.78%  '
Set-TextValue "E32" '  +4.55%  '
Set-TextValue "D33" '7.37'
Set-TextValue "E33" '  -3.05%  '
Set-TextValue "E34" '  -3.14%  '
Set-TextValue "D35" '40.47'
Set-TextValue "E35" '  -0.74%  '
Set-TextValue "D36" '0.999'
Set-TextValue "E36" '  -0.06%  '
Set-TextValue "D37" '56.78'
Set-TextValue "E37" '  -3.09%  '
Set-TextValue "D38" '0.0494'
Set-TextValue "E38" '  +0.50%  '
Set-TextValue "D39" '0.0₃0737'
Set-TextValue "E39" '  +32.21%  '
Set-TextValue "E40" '  +9.83%  '
Set-TextValue "E41" '  -0.09%  '
Set-TextValue "D42" '3.03'
Set-TextValue "E42" '  -2.07%  '
Set-TextValue "D43" '149.29'
Set-TextValue "E43" '  +1.38%  '
Set-TextValue "E44" '  +2.01%  '
Set-TextValue "D45" '3.28'
Set-TextValue "E45" '  -1.88%  '
Set-TextValue "D46" '0.315'
Set-TextValue "E46" '  -1.98%  '
Set-TextValue "D47" '4.32'
Set-TextValue "E47" '  -1.21%  '
Set-TextValue "E48" '  -4.77%  '
Set-TextValue "D49" '2.30'
Set-TextValue "E49" '  -4.94%  '

# Rows 50 and 51 swapped position in the refreshed ranking (Celestia dropped
# out of the top spot, BitcoinSV moved up) and both picked up new
# price/volume figures.
Set-TextValue "B50" 'BitcoinSV'
Set-TextValue "C50" 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue "D50" '115.63'
Set-TextValue "E50" '  +6.09%  '

Set-TextValue "B51" 'Celestia'
Set-TextValue "C51" 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue "D51" '15.65'
Set-TextValue "E51" '  -4.65%  '
